$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- logo_platform (sheet1): update logo URLs, add new "npai" row ---

# Set the new row's label first so its shared string is allocated before the
# relocated iaPos URL (keeps shared-string ordering aligned with the target).
$ws1.Range("A4").Value = "npai"

# Update the iaPos hyperlink display text / target to the new host.
$ws1.Range("B2").Value = "https://simcc.uesc.br/api/iapos.png"

# New row 4: npai platform logo.
$ws1.Range("B4").Value = "https://simcc.uesc.br/api/npai.png"

# Refresh the hyperlink relationships so B2/B3/B4 point at the right targets.
# (Deleting via a scoped range clears every hyperlink on the sheet in this
# engine, so we always rebuild the full set in ref order.)
$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://simcc.uesc.br/api/iapos.png")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://conectee.eng.ufmg.br/powerbi/conectee.png")
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://simcc.uesc.br/api/npai.png")

# Hyperlinks.Add stamps direct hyperlink formatting; re-apply the plain
# "Hiperlink" cell style so B2/B3/B4 match the original formatting.
$ws1.Range("B2").Style = "Hiperlink"
$ws1.Range("B3").Style = "Hiperlink"
$ws1.Range("B4").Style = "Hiperlink"

# --- image_painel (sheet2): no content changes ---

# --- Selection / active sheet: logo_platform, cell B4 ---
$ws1.Activate()
$ws1.Range("B4").Select()
